$d = $word.ActiveDocument

# Delete the four list-item paragraphs that were removed:
#   "Do the last acceptance tests and change in test report and software dev doc"
#   "Finish Release history - write down features in pre release"
#   "Enough bug fixes?"
#   "Add burn down chart in Software Development Document"
# They are four consecutive paragraphs right after the blank paragraph
# following the title, and right before the first "Go through the content
# in Software Development Document" paragraph.
$pStart = $d.Paragraphs.Item(3)
$pEnd = $d.Paragraphs.Item(6)
$delRange = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$delRange.Delete()

# Move the "_GoBack" bookmark from the (now trailing, blank) paragraph it
# currently lives in to the very start of the first remaining
# "Go through the content in Software Development Document" paragraph.
# Adding a bookmark named "_GoBack" moves/replaces the existing one.
$target = $d.Paragraphs.Item(3)
$insertionPoint = $d.Range($target.Range.Start, $target.Range.Start)
$d.Bookmarks.Add("_GoBack", $insertionPoint)
